$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (pushes "fossil_routes" and everything
# below it down by one row) to make room for the new
# "chemical_recycling_pyrolysis" parameter, right after
# "chemical_recycling_gasification" (row 9).
$ws.Rows.Item(10).Insert()

# Fill in the values for the newly inserted row.
$ws.Cells.Item(10, 1).Value = "chemical_recycling_pyrolysis"
$ws.Cells.Item(10, 2).Value = $true
